$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.1611213333333333
$ws.Cells.Item(2, 8).Value = 0.483364
$ws.Cells.Item(2, 9).Value = 0.001886845364621106
$ws.Cells.Item(2, 10).Value = 0.001927375876930126
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.104464
$ws.Cells.Item(2, 14).Value = 0.313392
$ws.Cells.Item(2, 15).Value = 0.02386241830158671
$ws.Cells.Item(2, 16).Value = 0.02789865426061783
$ws.Cells.Item(2, 17).Value = 0.01683137896533334
$ws.Cells.Item(2, 18).Value = 0.151482410688
$ws.Cells.Item(2, 19).Value = 0.00004502469336099872
$ws.Cells.Item(2, 20).Value = 0.00005377119322072868
$ws.Cells.Item(3, 7).Value = 0.1611213333333333
$ws.Cells.Item(3, 8).Value = 0.483364
$ws.Cells.Item(3, 9).Value = 0.001886845364621106
$ws.Cells.Item(3, 10).Value = 0.001927375876930126
$ws.Cells.Item(3, 15).Value = 0.02977487548338525
$ws.Cells.Item(3, 16).Value = 0.03481118075566868
$ws.Cells.Item(3, 17).Value = 0.02100173614311111
$ws.Cells.Item(3, 18).Value = 0.189015625288
$ws.Cells.Item(3, 19).Value = 0.00005618058578799607
$ws.Cells.Item(3, 20).Value = 0.00006709423003593003
$ws.Cells.Item(4, 7).Value = 0.1611213333333333
$ws.Cells.Item(4, 8).Value = 0.483364
$ws.Cells.Item(4, 9).Value = 0.001886845364621106
$ws.Cells.Item(4, 10).Value = 0.001927375876930126
$ws.Cells.Item(4, 13).Value = 1.913069666666667
$ws.Cells.Item(4, 14).Value = 5.739209
$ws.Cells.Item(4, 15).Value = 0.4369971341905063
$ws.Cells.Item(4, 16).Value = 0.5109135128542726
$ws.Cells.Item(4, 17).Value = 0.3082363354528889
$ws.Cells.Item(4, 18).Value = 2.774127019076
$ws.Cells.Item(4, 19).Value = 0.0008245460170000641
$ws.Cells.Item(4, 20).Value = 0.0009847223798729546
$ws.Cells.Item(5, 7).Value = 0.1611213333333333
$ws.Cells.Item(5, 8).Value = 0.483364
$ws.Cells.Item(5, 9).Value = 0.001886845364621106
$ws.Cells.Item(5, 10).Value = 0.001927375876930126
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.9000575
$ws.Cells.Item(5, 14).Value = 3.800115
$ws.Cells.Item(5, 15).Value = 0.4340248014824925
$ws.Cells.Item(5, 16).Value = 0.3382922810269175
$ws.Cells.Item(5, 17).Value = 0.30613979781
$ws.Cells.Item(5, 18).Value = 1.83683878686
$ws.Cells.Item(5, 19).Value = 0.0008189376848078364
$ws.Cells.Item(5, 20).Value = 0.0006520163818029477
$ws.Cells.Item(6, 7).Value = 0.1611213333333333
$ws.Cells.Item(6, 8).Value = 0.483364
$ws.Cells.Item(6, 9).Value = 0.001886845364621106
$ws.Cells.Item(6, 10).Value = 0.001927375876930126
$ws.Cells.Item(6, 13).Value = 0.329824
$ws.Cells.Item(6, 14).Value = 0.989472
$ws.Cells.Item(6, 15).Value = 0.07534077054202917
$ws.Cells.Item(6, 16).Value = 0.0880843711025235
$ws.Cells.Item(6, 17).Value = 0.05314168264533334
$ws.Cells.Item(6, 18).Value = 0.478275143808
$ws.Cells.Item(6, 19).Value = 0.0001421563836642101
$ws.Cells.Item(6, 20).Value = 0.0001697716919975649
$ws.Cells.Item(7, 9).Value = 0.07770616886214393
$ws.Cells.Item(7, 10).Value = 0.07937534159490063
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.104464
$ws.Cells.Item(7, 14).Value = 0.313392
$ws.Cells.Item(7, 15).Value = 0.02386241830158671
$ws.Cells.Item(7, 16).Value = 0.02789865426061783
$ws.Cells.Item(7, 17).Value = 0.6931686086133334
$ws.Cells.Item(7, 18).Value = 6.23851747752
$ws.Cells.Item(7, 19).Value = 0.001854257106002211
$ws.Cells.Item(7, 20).Value = 0.00221446521197457
$ws.Cells.Item(8, 9).Value = 0.07770616886214393
$ws.Cells.Item(8, 10).Value = 0.07937534159490063
$ws.Cells.Item(8, 15).Value = 0.02977487548338525
$ws.Cells.Item(8, 16).Value = 0.03481118075566868
$ws.Cells.Item(8, 19).Value = 0.002313691502161244
$ws.Cells.Item(8, 20).Value = 0.002763149363803032
$ws.Cells.Item(9, 9).Value = 0.07770616886214393
$ws.Cells.Item(9, 10).Value = 0.07937534159490063
$ws.Cells.Item(9, 13).Value = 1.913069666666667
$ws.Cells.Item(9, 14).Value = 5.739209
$ws.Cells.Item(9, 15).Value = 0.4369971341905063
$ws.Cells.Item(9, 16).Value = 0.5109135128542726
$ws.Cells.Item(9, 17).Value = 12.69413232332389
$ws.Cells.Item(9, 18).Value = 114.247190909915
$ws.Cells.Item(9, 19).Value = 0.03395737310168045
$ws.Cells.Item(9, 20).Value = 0.04055393460825854
$ws.Cells.Item(10, 9).Value = 0.07770616886214393
$ws.Cells.Item(10, 10).Value = 0.07937534159490063
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.9000575
$ws.Cells.Item(10, 14).Value = 3.800115
$ws.Cells.Item(10, 15).Value = 0.4340248014824925
$ws.Cells.Item(10, 16).Value = 0.3382922810269175
$ws.Cells.Item(10, 17).Value = 12.6077903733375
$ws.Cells.Item(10, 18).Value = 75.64674224002501
$ws.Cells.Item(10, 19).Value = 0.03372640451435706
$ws.Cells.Item(10, 20).Value = 0.0268520653654297
$ws.Cells.Item(11, 9).Value = 0.07770616886214393
$ws.Cells.Item(11, 10).Value = 0.07937534159490063
$ws.Cells.Item(11, 13).Value = 0.329824
$ws.Cells.Item(11, 14).Value = 0.989472
$ws.Cells.Item(11, 15).Value = 0.07534077054202917
$ws.Cells.Item(11, 16).Value = 0.0880843711025235
$ws.Cells.Item(11, 17).Value = 2.188540005813334
$ws.Cells.Item(11, 18).Value = 19.69686005232
$ws.Cells.Item(11, 19).Value = 0.005854442637942959
$ws.Cells.Item(11, 20).Value = 0.006991727045434796
$ws.Cells.Item(12, 7).Value = 40.78183766666667
$ws.Cells.Item(12, 8).Value = 122.345513
$ws.Cells.Item(12, 9).Value = 0.4775843134495767
$ws.Cells.Item(12, 10).Value = 0.4878430963142499
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.104464
$ws.Cells.Item(12, 14).Value = 0.313392
$ws.Cells.Item(12, 15).Value = 0.02386241830158671
$ws.Cells.Item(12, 16).Value = 0.02789865426061783
$ws.Cells.Item(12, 17).Value = 4.260233890010666
$ws.Cells.Item(12, 18).Value = 38.342105010096
$ws.Cells.Item(12, 19).Value = 0.0113963166618099
$ws.Cells.Item(12, 20).Value = 0.01361016587750054
$ws.Cells.Item(13, 7).Value = 40.78183766666667
$ws.Cells.Item(13, 8).Value = 122.345513
$ws.Cells.Item(13, 9).Value = 0.4775843134495767
$ws.Cells.Item(13, 10).Value = 0.4878430963142499
$ws.Cells.Item(13, 15).Value = 0.02977487548338525
$ws.Cells.Item(13, 16).Value = 0.03481118075566868
$ws.Cells.Item(13, 17).Value = 5.315803788282889
$ws.Cells.Item(13, 18).Value = 47.842234094546
$ws.Cells.Item(13, 19).Value = 0.01422001346577918
$ws.Cells.Item(13, 20).Value = 0.01698239420620044
$ws.Cells.Item(14, 7).Value = 40.78183766666667
$ws.Cells.Item(14, 8).Value = 122.345513
$ws.Cells.Item(14, 9).Value = 0.4775843134495767
$ws.Cells.Item(14, 10).Value = 0.4878430963142499
$ws.Cells.Item(14, 13).Value = 1.913069666666667
$ws.Cells.Item(14, 14).Value = 5.739209
$ws.Cells.Item(14, 15).Value = 0.4369971341905063
$ws.Cells.Item(14, 16).Value = 0.5109135128542726
$ws.Cells.Item(14, 17).Value = 78.01849659102412
$ws.Cells.Item(14, 18).Value = 702.166469319217
$ws.Cells.Item(14, 19).Value = 0.2087029763118055
$ws.Cells.Item(14, 20).Value = 0.2492456300596186
$ws.Cells.Item(15, 7).Value = 40.78183766666667
$ws.Cells.Item(15, 8).Value = 122.345513
$ws.Cells.Item(15, 9).Value = 0.4775843134495767
$ws.Cells.Item(15, 10).Value = 0.4878430963142499
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.9000575
$ws.Cells.Item(15, 14).Value = 3.800115
$ws.Cells.Item(15, 15).Value = 0.4340248014824925
$ws.Cells.Item(15, 16).Value = 0.3382922810269175
$ws.Cells.Item(15, 17).Value = 77.4878365223325
$ws.Cells.Item(15, 18).Value = 464.927019133995
$ws.Cells.Item(15, 19).Value = 0.207283436836105
$ws.Cells.Item(15, 20).Value = 0.1650335538353818
$ws.Cells.Item(16, 7).Value = 40.78183766666667
$ws.Cells.Item(16, 8).Value = 122.345513
$ws.Cells.Item(16, 9).Value = 0.4775843134495767
$ws.Cells.Item(16, 10).Value = 0.4878430963142499
$ws.Cells.Item(16, 13).Value = 0.329824
$ws.Cells.Item(16, 14).Value = 0.989472
$ws.Cells.Item(16, 15).Value = 0.07534077054202917
$ws.Cells.Item(16, 16).Value = 0.0880843711025235
$ws.Cells.Item(16, 17).Value = 13.45082882657067
$ws.Cells.Item(16, 18).Value = 121.057459439136
$ws.Cells.Item(16, 19).Value = 0.0359815701740771
$ws.Cells.Item(16, 20).Value = 0.0429713523355485
$ws.Cells.Item(17, 7).Value = 5.387083000000001
$ws.Cells.Item(17, 8).Value = 10.774166
$ws.Cells.Item(17, 9).Value = 0.06308657194606442
$ws.Cells.Item(17, 10).Value = 0.04296113827765565
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.104464
$ws.Cells.Item(17, 14).Value = 0.313392
$ws.Cells.Item(17, 15).Value = 0.02386241830158671
$ws.Cells.Item(17, 16).Value = 0.02789865426061783
$ws.Cells.Item(17, 17).Value = 0.5627562385120001
$ws.Cells.Item(17, 18).Value = 3.376537431072
$ws.Cells.Item(17, 19).Value = 0.001505398168990134
$ws.Cells.Item(17, 20).Value = 0.00119855794345091
$ws.Cells.Item(18, 7).Value = 5.387083000000001
$ws.Cells.Item(18, 8).Value = 10.774166
$ws.Cells.Item(18, 9).Value = 0.06308657194606442
$ws.Cells.Item(18, 10).Value = 0.04296113827765565
$ws.Cells.Item(18, 15).Value = 0.02977487548338525
$ws.Cells.Item(18, 16).Value = 0.03481118075566868
$ws.Cells.Item(18, 17).Value = 0.7021919034953334
$ws.Cells.Item(18, 18).Value = 4.213151420972
$ws.Cells.Item(18, 19).Value = 0.001878394824367693
$ws.Cells.Item(18, 20).Value = 0.001495527950052747
$ws.Cells.Item(19, 7).Value = 5.387083000000001
$ws.Cells.Item(19, 8).Value = 10.774166
$ws.Cells.Item(19, 9).Value = 0.06308657194606442
$ws.Cells.Item(19, 10).Value = 0.04296113827765565
$ws.Cells.Item(19, 13).Value = 1.913069666666667
$ws.Cells.Item(19, 14).Value = 5.739209
$ws.Cells.Item(19, 15).Value = 0.4369971341905063
$ws.Cells.Item(19, 16).Value = 0.5109135128542726
$ws.Cells.Item(19, 17).Value = 10.30586507911567
$ws.Cells.Item(19, 18).Value = 61.835190474694
$ws.Cells.Item(19, 19).Value = 0.02756865114633334
$ws.Cells.Item(19, 20).Value = 0.0219494260736552
$ws.Cells.Item(20, 7).Value = 5.387083000000001
$ws.Cells.Item(20, 8).Value = 10.774166
$ws.Cells.Item(20, 9).Value = 0.06308657194606442
$ws.Cells.Item(20, 10).Value = 0.04296113827765565
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 1.9000575
$ws.Cells.Item(20, 14).Value = 3.800115
$ws.Cells.Item(20, 15).Value = 0.4340248014824925
$ws.Cells.Item(20, 16).Value = 0.3382922810269175
$ws.Cells.Item(20, 17).Value = 10.2357674572725
$ws.Cells.Item(20, 18).Value = 40.94306982909
$ws.Cells.Item(20, 19).Value = 0.02738113686510158
$ws.Cells.Item(20, 20).Value = 0.01453342146346095
$ws.Cells.Item(21, 7).Value = 5.387083000000001
$ws.Cells.Item(21, 8).Value = 10.774166
$ws.Cells.Item(21, 9).Value = 0.06308657194606442
$ws.Cells.Item(21, 10).Value = 0.04296113827765565
$ws.Cells.Item(21, 13).Value = 0.329824
$ws.Cells.Item(21, 14).Value = 0.989472
$ws.Cells.Item(21, 15).Value = 0.07534077054202917
$ws.Cells.Item(21, 16).Value = 0.0880843711025235
$ws.Cells.Item(21, 17).Value = 1.776789263392
$ws.Cells.Item(21, 18).Value = 10.660735580352
$ws.Cells.Item(21, 19).Value = 0.004752990941271654
$ws.Cells.Item(21, 20).Value = 0.003784204847035847
$ws.Cells.Item(22, 7).Value = 32.42639166666667
$ws.Cells.Item(22, 8).Value = 97.279175
$ws.Cells.Item(22, 9).Value = 0.379736100377594
$ws.Cells.Item(22, 10).Value = 0.3878930479362637
$ws.Cells.Item(22, 11).Value = 2
$ws.Cells.Item(22, 12).Value = 0.6666666666666666
$ws.Cells.Item(22, 13).Value = 0.104464
$ws.Cells.Item(22, 14).Value = 0.313392
$ws.Cells.Item(22, 15).Value = 0.02386241830158671
$ws.Cells.Item(22, 16).Value = 0.02789865426061783
$ws.Cells.Item(22, 17).Value = 3.387390579066667
$ws.Cells.Item(22, 18).Value = 30.4865152116
$ws.Cells.Item(22, 19).Value = 0.009061421671423467
$ws.Cells.Item(22, 20).Value = 0.01082169403447108
$ws.Cells.Item(23, 7).Value = 32.42639166666667
$ws.Cells.Item(23, 8).Value = 97.279175
$ws.Cells.Item(23, 9).Value = 0.379736100377594
$ws.Cells.Item(23, 10).Value = 0.3878930479362637
$ws.Cells.Item(23, 15).Value = 0.02977487548338525
$ws.Cells.Item(23, 16).Value = 0.03481118075566868
$ws.Cells.Item(23, 17).Value = 4.226693683372223
$ws.Cells.Item(23, 18).Value = 38.04024315035
$ws.Cells.Item(23, 19).Value = 0.01130659510528914
$ws.Cells.Item(23, 20).Value = 0.01350301500557653
$ws.Cells.Item(24, 7).Value = 32.42639166666667
$ws.Cells.Item(24, 8).Value = 97.279175
$ws.Cells.Item(24, 9).Value = 0.379736100377594
$ws.Cells.Item(24, 10).Value = 0.3878930479362637
$ws.Cells.Item(24, 13).Value = 1.913069666666667
$ws.Cells.Item(24, 14).Value = 5.739209
$ws.Cells.Item(24, 15).Value = 0.4369971341905063
$ws.Cells.Item(24, 16).Value = 0.5109135128542726
$ws.Cells.Item(24, 17).Value = 62.03394629695278
$ws.Cells.Item(24, 18).Value = 558.305516672575
$ws.Cells.Item(24, 19).Value = 0.165943587613687
$ws.Cells.Item(24, 20).Value = 0.1981797997328672
$ws.Cells.Item(25, 7).Value = 32.42639166666667
$ws.Cells.Item(25, 8).Value = 97.279175
$ws.Cells.Item(25, 9).Value = 0.379736100377594
$ws.Cells.Item(25, 10).Value = 0.3878930479362637
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 1.9000575
$ws.Cells.Item(25, 14).Value = 3.800115
$ws.Cells.Item(25, 15).Value = 0.4340248014824925
$ws.Cells.Item(25, 16).Value = 0.3382922810269175
$ws.Cells.Item(25, 17).Value = 61.6120086841875
$ws.Cells.Item(25, 18).Value = 369.672052105125
$ws.Cells.Item(25, 19).Value = 0.164814885582121
$ws.Cells.Item(25, 20).Value = 0.1312212239808421
$ws.Cells.Item(26, 7).Value = 32.42639166666667
$ws.Cells.Item(26, 8).Value = 97.279175
$ws.Cells.Item(26, 9).Value = 0.379736100377594
$ws.Cells.Item(26, 10).Value = 0.3878930479362637
$ws.Cells.Item(26, 13).Value = 0.329824
$ws.Cells.Item(26, 14).Value = 0.989472
$ws.Cells.Item(26, 15).Value = 0.07534077054202917
$ws.Cells.Item(26, 16).Value = 0.0880843711025235
$ws.Cells.Item(26, 17).Value = 10.69500220506667
$ws.Cells.Item(26, 18).Value = 96.2550198456
$ws.Cells.Item(26, 19).Value = 0.02860961040507326
$ws.Cells.Item(26, 20).Value = 0.03416731518250678

Write-Output "Applied 298 cell updates"